$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 301-302 (Excel-style InsertRows), pushing the
# existing rows 301-318 down to 303-320 and carrying their formatting.
$ws.Rows("301:302").Insert()

# New row 301: Choclo / Dulce o Americano imported from Argentina
$ws.Range("A301").Value = 9
$ws.Range("B301").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C301").Value = "Metropolitana"
$ws.Range("D301").Value = 44516
$ws.Range("E301").Value = 13
$ws.Range("F301").Value = 100112024
$ws.Range("G301").Value = "Choclo"
$ws.Range("H301").Value = "Dulce o Americano"
$ws.Range("I301").Value = "Primera"
$ws.Range("J301").Value = 25
$ws.Range("K301").Value = 20000
$ws.Range("L301").Value = 23000
$ws.Range("M301").Value = 21440
$ws.Range("N301").Value = "$/caja 50 unidades"
$ws.Range("O301").Value = "Argentina"
$ws.Range("P301").Value = 429
$ws.Range("Q301").Value = 50
$ws.Range("R301").Value = "Hortaliza"

# New row 302: Choclo / Dulce o Americano from Provincia de Limari
$ws.Range("A302").Value = 9
$ws.Range("B302").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C302").Value = "Metropolitana"
$ws.Range("D302").Value = 44516
$ws.Range("E302").Value = 13
$ws.Range("F302").Value = 100112024
$ws.Range("G302").Value = "Choclo"
$ws.Range("H302").Value = "Dulce o Americano"
$ws.Range("I302").Value = "Primera"
$ws.Range("J302").Value = 34
$ws.Range("K302").Value = 30000
$ws.Range("L302").Value = 32000
$ws.Range("M302").Value = 31000
$ws.Range("N302").Value = "$/malla 60 unidades"
$ws.Range("O302").Value = "Provincia de Limarí"
$ws.Range("P302").Value = 517
$ws.Range("Q302").Value = 60
$ws.Range("R302").Value = "Hortaliza"
